$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set header cells I1 and J1, copying style from H1 (bold/border/center)
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for I and J columns, rows 2-81: row, I-value, J-value
$data = @(
    @(2, 9, 9),
    @(3, 8, 8),
    @(4, 7, 8),
    @(5, 5, 5),
    @(6, 7, 8),
    @(7, 7, 7),
    @(8, 8, 8),
    @(9, 8, 8),
    @(10, 9, 9),
    @(11, 9, 9),
    @(12, 5, 6),
    @(13, 7, 8),
    @(14, 9, 9),
    @(15, 7, 7),
    @(16, 6, 7),
    @(17, 8, 8),
    @(18, 7, 7),
    @(19, 7, 7),
    @(20, 9, 9),
    @(21, 10, 10),
    @(22, 6, 7),
    @(23, 8, 8),
    @(24, 9, 9),
    @(25, 5, 6),
    @(26, 8, 8),
    @(27, 9, 9),
    @(28, 8, 8),
    @(29, 7, 7),
    @(30, 8, 9),
    @(31, 5, 6),
    @(32, 7, 8),
    @(33, 6, 6),
    @(34, 6, 6),
    @(35, 8, 8),
    @(36, 7, 7),
    @(37, 8, 8),
    @(38, 8, 8),
    @(39, 6, 7),
    @(40, 6, 7),
    @(41, 6, 6),
    @(42, 4, 5),
    @(43, 6, 6),
    @(44, 6, 6),
    @(45, 9, 9),
    @(46, 8, 8),
    @(47, 6, 6),
    @(48, 5, 5),
    @(49, 12, 12),
    @(50, 3, 4),
    @(51, 7, 8),
    @(52, 6, 6),
    @(53, 11, 12),
    @(54, 7, 7),
    @(55, 5, 5),
    @(56, 6, 7),
    @(57, 7, 7),
    @(58, 8, 8),
    @(59, 8, 9),
    @(60, 9, 9),
    @(61, 8, 8),
    @(62, 5, 5),
    @(63, 6, 6),
    @(64, 8, 8),
    @(65, 6, 7),
    @(66, 8, 9),
    @(67, 6, 6),
    @(68, 5, 6),
    @(69, 9, 9),
    @(70, 8, 8),
    @(71, 8, 8),
    @(72, 5, 5),
    @(73, 9, 9),
    @(74, 8, 9),
    @(75, 4, 4),
    @(76, 6, 6),
    @(77, 6, 6),
    @(78, 6, 6),
    @(79, 9, 9),
    @(80, 6, 6),
    @(81, 5, 5)
)

foreach ($item in $data) {
    $row = $item[0]
    $iVal = $item[1]
    $jVal = $item[2]
    $ws.Cells.Item($row, 9).Value = $iVal
    $ws.Cells.Item($row, 10).Value = $jVal
}

Write-Output "Applied I0/IF columns"